$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Description" column header to include the note about leaving
# a single blank space when there is no description.
$ws.Range("C1").Value = "Description `n(설명이 없으면 빈칸하나 입력)"

# Move the active selection to C18 (first empty row below the table),
# matching where the author left off after editing the header.
$ws.Range("C18").Select()
